# Removed duplicated row in config files.
# Row 19 on the "Workflow" sheet is an exact duplicate of row 18
# (same "Open Browser not being used" check). Delete it; Excel will
# shift the following rows up and auto-adjust the sheet dimension and
# the data-validation ranges that reference the sheet's row extent.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Workflow")

$ws.Rows("19:19").Delete()
